$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.944.35'
$ws.Range('E2').Value = '  +0.68%  '
$ws.Range('D3').Value = '1.748.88'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9962'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9970'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5182'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2827'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +8.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.68'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06137'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').Value = '1.740.08'
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06999'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.50'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6447'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.534'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '77.10'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9954'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9969'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').Value = '25.967.17'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('E20').Value = '  -0.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006647'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.26%  '
$ws.Range('D22').Value = '1.962.52'
$ws.Range('E22').Value = '  -0.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.138'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.587'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.168'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.23'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.508'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.35%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.832'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.59%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.12'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '103.27'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08319'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.645'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.60%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.445'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04419'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.612'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9878'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6121'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.682'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01578'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.949'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9954'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '100.77'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3877'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7311'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.967'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.24%  '
$ws.Range('E46').Value = '  -0.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.401'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +7.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1116'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '52.73'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '29.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.61%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3428'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.30%  '
